$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the "purpose" column (E2:E13) to the new value "fullRNASEQ"
$ws.Range("E2:E13").Value = "fullRNASEQ"

# Update the active selection to match the saved state (E13 only)
$ws.Range("E13").Select()
